$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even when it looks like a number
# (e.g. "243.94", "14.25"), so Excel does not silently coerce it to a
# numeric cell -- matches the source data which stores these as strings.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "36.627.42"
$ws.Range("E2").Value = "  -0.06%  "
Set-TextValue $ws.Range("D3") "1.972.97"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "243.94"
$ws.Range("E5").Value = "  -0.21%  "
Set-TextValue $ws.Range("D6") "0.627"
$ws.Range("E6").Value = "  +2.12%  "
Set-TextValue $ws.Range("D7") "60.15"
$ws.Range("E7").Value = "  +3.60%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("E11").Value = "  +0.75%  "
Set-TextValue $ws.Range("D12") "14.25"
$ws.Range("E12").Value = "  +4.41%  "
Set-TextValue $ws.Range("D13") "0.842"
$ws.Range("E13").Value = "  +2.91%  "
Set-TextValue $ws.Range("D14") "2.257.95"
$ws.Range("E14").Value = "  +0.22%  "
Set-TextValue $ws.Range("D15") "21.70"
$ws.Range("E15").Value = "  -1.34%  "
Set-TextValue $ws.Range("D16") "5.35"
$ws.Range("E16").Value = "  +1.91%  "
Set-TextValue $ws.Range("D17") "1.966.19"
$ws.Range("E17").Value = "  +0.35%  "
Set-TextValue $ws.Range("D18") "36.552.14"
$ws.Range("E18").Value = "  -0.17%  "
Set-TextValue $ws.Range("D19") "69.78"
$ws.Range("E19").Value = "  +0.26%  "
Set-TextValue $ws.Range("D20") "0.0₃0855"
$ws.Range("E20").Value = "  -0.49%  "
Set-TextValue $ws.Range("D21") "229.57"
$ws.Range("E21").Value = "  +0.84%  "
Set-TextValue $ws.Range("D22") "5.09"
$ws.Range("E22").Value = "  +0.17%  "
Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("E25").Value = "  +0.55%  "
Set-TextValue $ws.Range("D26") "0.145"
$ws.Range("E26").Value = "  +7.34%  "
Set-TextValue $ws.Range("D27") "9.13"
$ws.Range("E27").Value = "  -2.05%  "
Set-TextValue $ws.Range("D28") "162.57"
$ws.Range("E28").Value = "  +1.14%  "
Set-TextValue $ws.Range("D29") "19.40"
$ws.Range("E29").Value = "  +0.26%  "
Set-TextValue $ws.Range("D30") "1.33"
$ws.Range("E30").Value = "  +19.05%  "
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("E32").Value = "  +2.93%  "
Set-TextValue $ws.Range("D33") "0.0615"
$ws.Range("E33").Value = "  -0.22%  "
Set-TextValue $ws.Range("D34") "4.50"
$ws.Range("E34").Value = "  +6.57%  "
Set-TextValue $ws.Range("D35") "2.29"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("E36").Value = "  +0.14%  "
Set-TextValue $ws.Range("D37") "3.35"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("E39").Value = "  -13.62%  "
Set-TextValue $ws.Range("D40") "0.0969"
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("E43").Value = "  -0.88%  "
Set-TextValue $ws.Range("D44") "15.93"
$ws.Range("E44").Value = "  -0.21%  "
Set-TextValue $ws.Range("D45") "1.365.95"
$ws.Range("E45").Value = "  +1.13%  "
Set-TextValue $ws.Range("D46") "89.18"
$ws.Range("E46").Value = "  +2.20%  "
Set-TextValue $ws.Range("D47") "1.03"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("E49").Value = "  -0.65%  "
Set-TextValue $ws.Range("D50") "46.14"
$ws.Range("E50").Value = "  +6.54%  "
Set-TextValue $ws.Range("D51") "2.152.40"
$ws.Range("E51").Value = "  +0.39%  "
